# Mondelez DMI KPI template - ticket fixes (11141 / 1123 / 11104)
#
# The only real data change is on VTW_POINTS!C8 ("Pallet - Full - PRD" row):
# the "multiple" flag is flipped from "n" to "y". The rest of the diff
# (row stubs 42-54 below the table, the scrolled/selected view state, and
# a handful of sub-pixel column-width nudges) are left-over view/cosmetic
# artifacts from the range the author selected (A44:N54) while working on
# the sheet - they carry no workbook content and are reproduced here only
# to the extent the Excel object model actually exposes them.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("VTW_POINTS")
$ws2 = $wb.Worksheets.Item("GOLD_ZONE")

# --- the actual content edit -------------------------------------------------
$ws1.Range("C8").Value = "y"

# --- cosmetic: materialize the empty (height-only) rows 42-54 below the table
$ws1.Rows("42:54").RowHeight = 12.8

# --- cosmetic: column widths nudged slightly (best attainable approximation
#     given the Excel ColumnWidth <-> stored-width quantization) ------------
$ws1.Columns.Item(1).ColumnWidth = 26.5
$ws1.Columns.Item(2).ColumnWidth = 8.666666666666666

$ws2.Columns.Item(1).ColumnWidth = 15.833333333333334

# --- cosmetic: selection/scroll state ---------------------------------------
# GOLD_ZONE picked up "A44:N54" as a second selection area (its own A3
# selection stays primary); select it first so VTW_POINTS ends up as the
# active/tab-selected sheet, matching the final workbook state.
[void]$ws2.Activate()
[void]$ws2.Range("A3").Select()

[void]$ws1.Activate()
[void]$ws1.Range("A44:N54").Select()
